$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Estado de Cuenta" EC data table refresh (previous periods removed, new periods/workers added)
# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora, F = Valor Mora, G = Salario Basico
$ws.Cells.Item(16, 3).Value = '1102799537'
$ws.Cells.Item(16, 4).Value = 'JULIO GUSTAVO SILGADO LACAYO'
$ws.Cells.Item(16, 5).Value = '2112'
$ws.Cells.Item(16, 6).Value = 53333
$ws.Cells.Item(16, 7).Value = 2500000
$ws.Cells.Item(17, 3).Value = '1102799537'
$ws.Cells.Item(17, 4).Value = 'JULIO GUSTAVO SILGADO LACAYO'
$ws.Cells.Item(17, 5).Value = '2111'
$ws.Cells.Item(17, 6).Value = 100000
$ws.Cells.Item(17, 7).Value = 2500000
$ws.Cells.Item(18, 3).Value = '1102799537'
$ws.Cells.Item(18, 4).Value = 'JULIO GUSTAVO SILGADO LACAYO'
$ws.Cells.Item(18, 5).Value = '2110'
$ws.Cells.Item(18, 6).Value = 100000
$ws.Cells.Item(18, 7).Value = 2500000
$ws.Cells.Item(19, 3).Value = '1102799537'
$ws.Cells.Item(19, 4).Value = 'JULIO GUSTAVO SILGADO LACAYO'
$ws.Cells.Item(19, 5).Value = '2109'
$ws.Cells.Item(19, 6).Value = 100000
$ws.Cells.Item(19, 7).Value = 2500000
$ws.Cells.Item(20, 3).Value = '1102799537'
$ws.Cells.Item(20, 4).Value = 'JULIO GUSTAVO SILGADO LACAYO'
$ws.Cells.Item(20, 5).Value = '2108'
$ws.Cells.Item(20, 6).Value = 100000
$ws.Cells.Item(20, 7).Value = 2500000
$ws.Cells.Item(21, 3).Value = '1102799537'
$ws.Cells.Item(21, 4).Value = 'JULIO GUSTAVO SILGADO LACAYO'
$ws.Cells.Item(21, 5).Value = '2107'
$ws.Cells.Item(21, 6).Value = 100000
$ws.Cells.Item(21, 7).Value = 2500000
$ws.Cells.Item(22, 3).Value = '1102799537'
$ws.Cells.Item(22, 4).Value = 'JULIO GUSTAVO SILGADO LACAYO'
$ws.Cells.Item(22, 5).Value = '2106'
$ws.Cells.Item(22, 6).Value = 100000
$ws.Cells.Item(22, 7).Value = 2500000
$ws.Cells.Item(23, 3).Value = '9101392'
$ws.Cells.Item(23, 4).Value = 'ARBEY SAMIR BELLO LOZANO'
$ws.Cells.Item(23, 5).Value = '2112'
$ws.Cells.Item(23, 6).Value = 53333
$ws.Cells.Item(23, 7).Value = 4500000
$ws.Cells.Item(24, 3).Value = '9101392'
$ws.Cells.Item(24, 4).Value = 'ARBEY SAMIR BELLO LOZANO'
$ws.Cells.Item(24, 5).Value = '2111'
$ws.Cells.Item(24, 6).Value = 100000
$ws.Cells.Item(24, 7).Value = 4500000
$ws.Cells.Item(25, 3).Value = '9101392'
$ws.Cells.Item(25, 4).Value = 'ARBEY SAMIR BELLO LOZANO'
$ws.Cells.Item(25, 5).Value = '2110'
$ws.Cells.Item(25, 6).Value = 140000
$ws.Cells.Item(25, 7).Value = 4500000
$ws.Cells.Item(26, 3).Value = '9101392'
$ws.Cells.Item(26, 4).Value = 'ARBEY SAMIR BELLO LOZANO'
$ws.Cells.Item(26, 5).Value = '2109'
$ws.Cells.Item(26, 6).Value = 140000
$ws.Cells.Item(26, 7).Value = 4500000
$ws.Cells.Item(27, 3).Value = '9101392'
$ws.Cells.Item(27, 4).Value = 'ARBEY SAMIR BELLO LOZANO'
$ws.Cells.Item(27, 5).Value = '2108'
$ws.Cells.Item(27, 6).Value = 140000
$ws.Cells.Item(27, 7).Value = 4500000
$ws.Cells.Item(28, 3).Value = '9101392'
$ws.Cells.Item(28, 4).Value = 'ARBEY SAMIR BELLO LOZANO'
$ws.Cells.Item(28, 5).Value = '2107'
$ws.Cells.Item(28, 6).Value = 140000
$ws.Cells.Item(28, 7).Value = 4500000
$ws.Cells.Item(29, 3).Value = '9101392'
$ws.Cells.Item(29, 4).Value = 'ARBEY SAMIR BELLO LOZANO'
$ws.Cells.Item(29, 5).Value = '2106'
$ws.Cells.Item(29, 6).Value = 140000
$ws.Cells.Item(29, 7).Value = 4500000
$ws.Cells.Item(30, 3).Value = '1049937385'
$ws.Cells.Item(30, 4).Value = 'DIANA PATRICIA RIOS VILLERO'
$ws.Cells.Item(30, 5).Value = '2112'
$ws.Cells.Item(30, 6).Value = 23467
$ws.Cells.Item(30, 7).Value = 1100000
$ws.Cells.Item(31, 3).Value = '1049937385'
$ws.Cells.Item(31, 4).Value = 'DIANA PATRICIA RIOS VILLERO'
$ws.Cells.Item(31, 5).Value = '2111'
$ws.Cells.Item(31, 6).Value = 44000
$ws.Cells.Item(31, 7).Value = 1100000
$ws.Cells.Item(32, 3).Value = '1049937385'
$ws.Cells.Item(32, 4).Value = 'DIANA PATRICIA RIOS VILLERO'
$ws.Cells.Item(32, 5).Value = '2110'
$ws.Cells.Item(32, 6).Value = 44000
$ws.Cells.Item(32, 7).Value = 1100000
$ws.Cells.Item(33, 3).Value = '1049937385'
$ws.Cells.Item(33, 4).Value = 'DIANA PATRICIA RIOS VILLERO'
$ws.Cells.Item(33, 5).Value = '2109'
$ws.Cells.Item(33, 6).Value = 44000
$ws.Cells.Item(33, 7).Value = 1100000
$ws.Cells.Item(34, 3).Value = '1049937385'
$ws.Cells.Item(34, 4).Value = 'DIANA PATRICIA RIOS VILLERO'
$ws.Cells.Item(34, 5).Value = '2108'
$ws.Cells.Item(34, 6).Value = 44000
$ws.Cells.Item(34, 7).Value = 1100000
$ws.Cells.Item(35, 3).Value = '1049937385'
$ws.Cells.Item(35, 4).Value = 'DIANA PATRICIA RIOS VILLERO'
$ws.Cells.Item(35, 5).Value = '2107'
$ws.Cells.Item(35, 6).Value = 44000
$ws.Cells.Item(35, 7).Value = 1100000
$ws.Cells.Item(36, 3).Value = '1049937385'
$ws.Cells.Item(36, 4).Value = 'DIANA PATRICIA RIOS VILLERO'
$ws.Cells.Item(36, 5).Value = '2106'
$ws.Cells.Item(36, 6).Value = 44000
$ws.Cells.Item(36, 7).Value = 1100000
$ws.Cells.Item(37, 3).Value = '1049936487'
$ws.Cells.Item(37, 4).Value = 'ENEIDIS ESALAS URBINA'
$ws.Cells.Item(37, 5).Value = '2112'
$ws.Cells.Item(37, 6).Value = 23467
$ws.Cells.Item(37, 7).Value = 1400000
$ws.Cells.Item(38, 3).Value = '1049936487'
$ws.Cells.Item(38, 4).Value = 'ENEIDIS ESALAS URBINA'
$ws.Cells.Item(38, 5).Value = '2111'
$ws.Cells.Item(38, 6).Value = 44000
$ws.Cells.Item(38, 7).Value = 1400000
$ws.Cells.Item(39, 3).Value = '1049936487'
$ws.Cells.Item(39, 4).Value = 'ENEIDIS ESALAS URBINA'
$ws.Cells.Item(39, 5).Value = '2110'
$ws.Cells.Item(39, 6).Value = 44000
$ws.Cells.Item(39, 7).Value = 1400000
$ws.Cells.Item(40, 3).Value = '1049936487'
$ws.Cells.Item(40, 4).Value = 'ENEIDIS ESALAS URBINA'
$ws.Cells.Item(40, 5).Value = '2109'
$ws.Cells.Item(40, 6).Value = 44000
$ws.Cells.Item(40, 7).Value = 1400000
$ws.Cells.Item(41, 3).Value = '1049936487'
$ws.Cells.Item(41, 4).Value = 'ENEIDIS ESALAS URBINA'
$ws.Cells.Item(41, 5).Value = '2108'
$ws.Cells.Item(41, 6).Value = 44000
$ws.Cells.Item(41, 7).Value = 1400000
$ws.Cells.Item(42, 3).Value = '1049936487'
$ws.Cells.Item(42, 4).Value = 'ENEIDIS ESALAS URBINA'
$ws.Cells.Item(42, 5).Value = '2107'
$ws.Cells.Item(42, 6).Value = 44000
$ws.Cells.Item(42, 7).Value = 1400000
$ws.Cells.Item(43, 3).Value = '1049936487'
$ws.Cells.Item(43, 4).Value = 'ENEIDIS ESALAS URBINA'
$ws.Cells.Item(43, 5).Value = '2106'
$ws.Cells.Item(43, 6).Value = 17600
$ws.Cells.Item(43, 7).Value = 1400000
$ws.Cells.Item(44, 3).Value = '1049932362'
$ws.Cells.Item(44, 4).Value = 'WILBERTO PEREZ GUTIERREZ'
$ws.Cells.Item(44, 5).Value = '2112'
$ws.Cells.Item(44, 6).Value = 32000
$ws.Cells.Item(44, 7).Value = 1500000
$ws.Cells.Item(45, 3).Value = '1049932362'
$ws.Cells.Item(45, 4).Value = 'WILBERTO PEREZ GUTIERREZ'
$ws.Cells.Item(45, 5).Value = '2111'
$ws.Cells.Item(45, 6).Value = 60000
$ws.Cells.Item(45, 7).Value = 1500000
$ws.Cells.Item(46, 3).Value = '1049932362'
$ws.Cells.Item(46, 4).Value = 'WILBERTO PEREZ GUTIERREZ'
$ws.Cells.Item(46, 5).Value = '2110'
$ws.Cells.Item(46, 6).Value = 60000
$ws.Cells.Item(46, 7).Value = 1500000
$ws.Cells.Item(47, 3).Value = '1049932362'
$ws.Cells.Item(47, 4).Value = 'WILBERTO PEREZ GUTIERREZ'
$ws.Cells.Item(47, 5).Value = '2109'
$ws.Cells.Item(47, 6).Value = 60000
$ws.Cells.Item(47, 7).Value = 1500000
$ws.Cells.Item(48, 3).Value = '1049932362'
$ws.Cells.Item(48, 4).Value = 'WILBERTO PEREZ GUTIERREZ'
$ws.Cells.Item(48, 5).Value = '2108'
$ws.Cells.Item(48, 6).Value = 60000
$ws.Cells.Item(48, 7).Value = 1500000
$ws.Cells.Item(49, 3).Value = '1049932362'
$ws.Cells.Item(49, 4).Value = 'WILBERTO PEREZ GUTIERREZ'
$ws.Cells.Item(49, 5).Value = '2107'
$ws.Cells.Item(49, 6).Value = 60000
$ws.Cells.Item(49, 7).Value = 1500000
$ws.Cells.Item(50, 3).Value = '1049932362'
$ws.Cells.Item(50, 4).Value = 'WILBERTO PEREZ GUTIERREZ'
$ws.Cells.Item(50, 5).Value = '2106'
$ws.Cells.Item(50, 6).Value = 60000
$ws.Cells.Item(50, 7).Value = 1500000
